# Append the new ParamStudy rows (cases 66-73) below the existing data
# (previously ending at row 68) and update the sheet's active selection to
# the new last cell, matching the author's "param study general work" edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, CASE, ZETA) -- PTAU is the same constant (0.066) for every new row
$newRows = @(
    @(69, 66, 16),
    @(70, 66, 18),
    @(71, 66, 20),
    @(72, 67, 18),
    @(73, 67, 19),
    @(74, 67, 20),
    @(75, 68, 20),
    @(76, 68, 21),
    @(77, 68, 22),
    @(78, 69, 20),
    @(79, 69, 25),
    @(80, 69, 30),
    @(81, 70, 16),
    @(82, 70, 18),
    @(83, 70, 20),
    @(84, 71, 20),
    @(85, 71, 30),
    @(86, 71, 50),
    @(87, 72, 20),
    @(88, 72, 50),
    @(89, 72, 100),
    @(90, 73, 20),
    @(91, 73, 200),
    @(92, 73, 2000)
)

$ptau = 0.066
$firstRow = 69
$lastRow = 92

foreach ($entry in $newRows) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 5).Value = $ptau
}

# Copy the formatting (style index 1, used by every data row) from the
# last existing row onto the newly written B and E columns.
$ws.Range("B68").Copy()
$ws.Range("B$firstRow`:B$lastRow").PasteSpecial(-4122)
$ws.Range("E68").Copy()
$ws.Range("E$firstRow`:E$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to the new final cell, mirroring the
# scrolled/selected view recorded in the workbook after the edit.
$null = $ws.Range("A$lastRow").Select()
